$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (A1:D1) from descriptive Spanish labels to short codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/la/las/el/los/y) in
#    state (col A) and municipality (col B) names, e.g. "Pabellón de Arteaga"
#    -> "Pabellón De Arteaga". Built from a direct before/after cell map.
$cellUpdates = @{
    'B4' = 'Pabellón De Arteaga'
    'B5' = 'Rincón De Romos'
    'B24' = 'Amatenango De La Frontera'
    'B30' = 'Comitán De Domínguez'
    'B48' = 'Salto De Agua'
    'B64' = 'Hidalgo Del Parral'
    'B66' = 'San Francisco Del Oro'
    'B83' = 'Villa De Álvarez'
    'A85' = 'Ciudad De México'
    'B104' = 'Pánuco De Coronado'
    'B107' = 'San Luis Del Cordero'
    'A112' = 'Estado De México'
    'B112' = 'Acambay De Ruíz Castañeda'
    'B114' = 'Almoloya Del Río'
    'B123' = 'Ecatepec De Morelos'
    'B125' = 'Ixtapan De La Sal'
    'B132' = 'Naucalpan De Juárez'
    'B138' = 'San Felipe Del Progreso'
    'B140' = 'Soyaniquilpan De Juárez'
    'B150' = 'Tlalnepantla De Baz'
    'B153' = 'Villa Del Carbón'
    'B161' = 'Apaseo El Alto'
    'B167' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B178' = 'San Diego De La Unión'
    'B180' = 'San Francisco Del Rincón'
    'B182' = 'San Luis De La Paz'
    'B183' = 'Silao De La Victoria'
    'B186' = 'Valle De Santiago'
    'B189' = 'Acapulco De Juárez'
    'B191' = 'Ajuchitlán Del Progreso'
    'B192' = 'Alcozauca De Guerrero'
    'B195' = 'Atlamajalcingo Del Monte'
    'B197' = 'Atoyac De Álvarez'
    'B198' = 'Ayutla De Los Libres'
    'B201' = 'Chilapa De Álvarez'
    'B202' = 'Chilpancingo De Los Bravo'
    'B203' = 'Coahuayutla De José María Izazaga'
    'B207' = 'Coyuca De Benítez'
    'B208' = 'Coyuca De Catalán'
    'B212' = 'Cuetzala Del Progreso'
    'B216' = 'Iguala De La Independencia'
    'B218' = 'Ixcateopan De Cuauhtémoc'
    'B219' = 'Zihuatanejo De Azueta'
    'B220' = 'La Unión De Isidoro Montes De Oca'
    'B222' = 'Mártir De Cuilapan'
    'B233' = 'Taxco De Alarcón'
    'B235' = 'Técpan De Galeana'
    'B237' = 'Tepecoacuilco De Trujano'
    'B238' = 'Tixtla De Guerrero'
    'B242' = 'Tlalixtaquilla De Maldonado'
    'B243' = 'Tlapa De Comonfort'
    'B252' = 'Agua Blanca De Iturbide'
    'B256' = 'Atotonilco El Grande'
    'B258' = 'Cuautepec De Hinojosa'
    'B261' = 'Huasca De Ocampo'
    'B265' = 'Mineral De La Reforma'
    'B266' = 'Mineral Del Monte'
    'B267' = 'Nopala De Villagrán'
    'B268' = 'Pachuca De Soto'
    'B272' = 'Tenango De Doria'
    'B273' = 'Tepeji Del Río De Ocampo'
    'B274' = 'Tezontepec De Aldama'
    'B277' = 'Tula De Allende'
    'B278' = 'Tulancingo De Bravo'
    'B280' = 'Zacualtipán De Ángeles'
    'B283' = 'Acatlán De Juárez'
    'B284' = 'Ahualulco De Mercado'
    'B288' = 'Atemajac De Brizuela'
    'B290' = 'Atotonilco El Alto'
    'B292' = 'Autlán De Navarro'
    'B297' = 'Cañadas De Obregón'
    'B302' = 'Cuautitlán De García Barragán'
    'B311' = 'Huejuquilla El Alto'
    'B313' = 'Jilotlán De Los Dolores'
    'B315' = 'Lagos De Moreno'
    'B323' = 'San Diego De Alejandría'
    'B325' = 'San Juan De Los Lagos'
    'B327' = 'San Martín De Bolaños'
    'B329' = 'San Miguel El Alto'
    'B331' = 'Talpa De Allende'
    'B333' = 'Techaluta De Montenegro'
    'B336' = 'Teocuitatlán De Corona'
    'B337' = 'Tepatitlán De Morelos'
    'B339' = 'Tlajomulco De Zúñiga'
    'B345' = 'Unión De Tula'
    'B346' = 'Valle De Guadalupe'
    'B350' = 'Yahualica De González Gallo'
    'B351' = 'Zacoalco De Torres'
    'B354' = 'Zapotitlán De Vadillo'
    'B355' = 'Zapotlán Del Rey'
    'B356' = 'Zapotlán El Grande'
    'B380' = 'Coalcomán De Vázquez Pallares'
    'B382' = 'Cojumatlán De Régules'
    'B464' = 'Puente De Ixtla'
    'B468' = 'Tetela Del Volcán'
    'B469' = 'Tlaltizapán De Zapata'
    'B477' = 'Amatlán De Cañas'
    'B486' = 'Santa María Del Oro'
    'B497' = 'Acatlán De Pérez Figueroa'
    'B498' = 'Ayoquezco De Aldama'
    'B500' = 'Chalcatongo De Hidalgo'
    'B502' = 'Coicoyán De Las Flores'
    'B503' = 'Constancia Del Rosario'
    'B504' = 'Cuilápam De Guerrero'
    'B505' = 'Guadalupe De Ramírez'
    'B506' = 'Heroica Ciudad De Ejutla De Crespo'
    'B507' = 'Heroica Ciudad De Huajuapan De León'
    'B508' = 'Heroica Ciudad De Tlaxiaco'
    'B509' = 'Ixtlán De Juárez'
    'B510' = 'Heroica Ciudad De Juchitán De Zaragoza'
    'B517' = 'Mariscala De Juárez'
    'B519' = 'Mazatlán Villa De Flores'
    'B521' = 'Miahuatlán De Porfirio Díaz'
    'B522' = 'Mixistlán De La Reforma'
    'B524' = 'Nejapa De Madero'
    'B525' = 'Oaxaca De Juárez'
    'B526' = 'Ocotlán De Morelos'
    'B527' = 'Putla Villa De Guerrero'
    'B540' = 'San Antonino El Alto'
    'B558' = 'San José Del Peñasco'
    'B584' = 'San Miguel Del Puerto'
    'B585' = 'San Miguel El Grande'
    'B597' = 'San Pablo Villa De Mitla'
    'B600' = 'San Pedro El Alto'
    'B608' = 'San Pedro Y San Pablo Ayutla'
    'B609' = 'San Pedro Y San Pablo Tequixtepec'
    'B617' = 'Santa Ana Del Valle'
    'B628' = 'Santa Cruz Tacache De Mina'
    'B631' = 'Santa Inés Del Monte'
    'B684' = 'Santo Domingo De Morelos'
    'B696' = 'Tanetze De Zaragoza'
    'B697' = 'Teotitlán De Flores Magón'
    'B698' = 'Tepelmeme Villa De Morelos'
    'B699' = 'Tezoatlán De Segura Y Luna'
    'B700' = 'Tlacolula De Matamoros'
    'B702' = 'Tlalixtac De Cabrera'
    'B703' = 'Villa De Chilapa De Díaz'
    'B704' = 'Villa De Tututepec De Melchor Ocampo'
    'B705' = 'Villa De Zaachila'
    'B706' = 'Villa Sola De Vega'
    'B708' = 'Zimatlán De Álvarez'
    'B719' = 'Chalchicomula De Sesma'
    'B727' = 'Cuayuca De Andrade'
    'B728' = 'Cuetzalan Del Progreso'
    'B737' = 'Izúcar De Matamoros'
    'B740' = 'Mazapiltepec De Juárez'
    'B744' = 'Palmar De Bravo'
    'B755' = 'San Salvador El Verde'
    'B759' = 'Tecali De Herrera'
    'B765' = 'Tepanco De López'
    'B766' = 'Tepatlaxco De Hidalgo'
    'B768' = 'Tepexi De Rodríguez'
    'B769' = 'Tepeyahualco De Cuauhtémoc'
    'B771' = 'Tlacotepec De Benito Juárez'
    'B776' = 'Totoltepec De Guerrero'
    'B789' = 'Amealco De Bonfil'
    'B790' = 'Cadereyta De Montes'
    'B796' = 'San Juan Del Río'
    'B807' = 'Mexquitic De Carmona'
    'B815' = 'Villa De Ramos'
    'B841' = 'Nacozari De García'
    'B860' = 'Contla De Juan Cuamatzi'
    'B870' = 'Alto Lucero De Gutiérrez Barrios'
    'B871' = 'Amatlán De Los Reyes'
    'B879' = 'Castillo De Teayo'
    'B886' = 'Cosamaloapan De Carpio'
    'B898' = 'Hueyapan De Ocampo'
    'B901' = 'Ixhuatlán Del Café'
    'B910' = 'Lerdo De Tejada'
    'B914' = 'Martínez De La Torre'
    'B916' = 'Mixtla De Altamirano'
    'B925' = 'Poza Rica De Hidalgo'
    'B932' = 'Soledad De Doblado'
    'B950' = 'Zozocolco De Hidalgo'
    'B982' = 'Jiménez Del Teul'
    'B984' = 'Mezquital Del Oro'
    'B987' = 'Nochistlán De Mejía'
    'B991' = 'Teúl De González Ortega'
    'B992' = 'Tlaltenango De Sánchez Román'
    'B993' = 'Trinidad García De La Cadena'
    'B995' = 'Villa De Cos'
}

foreach ($ref in $cellUpdates.Keys) {
    $ws.Range($ref).Value = $cellUpdates[$ref]
}

# 3. Drop the trailing metadata/footnote rows (1001-1005): sample size,
#    source, author and date notes that no longer belong in the clean data.
$ws.Range("A1001:D1005").EntireRow.Delete() | Out-Null
